$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 221, shifting existing rows 221:267 down to 222:268.
$ws.Rows("221:221").Insert()

# Populate the newly inserted row 221 with the new weekly record.
$ws.Range("A221").Value = 10
$ws.Range("B221").Value = "Vega Modelo de Temuco"
$ws.Range("C221").Value = "La Araucanía"
$ws.Range("D221").Value = 44785
$ws.Range("E221").Value = 9
$ws.Range("F221").Value = 100112043
$ws.Range("G221").Value = "Pepino dulce"
$ws.Range("H221").Value = "Cultivar IV Región"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 80
$ws.Range("K221").Value = 18000
$ws.Range("L221").Value = 18000
$ws.Range("M221").Value = 18000
$ws.Range("N221").Value = "$/bandeja 18 kilos"
$ws.Range("O221").Value = "Provincia de Limarí"
$ws.Range("P221").Value = 1000
$ws.Range("Q221").Value = 18
$ws.Range("R221").Value = "Hortaliza"
